{"js": "// Update the multiplication-problem cells in the table with new operands.\n// Each \"old\" string is unique in the document, so a plain search + full\n// replace of the matched range keeps the existing run formatting\n// (font, size, etc.) intact while swapping only the visible text.\nconst pairs = [\n  [\"703\u00d76=\", \"442\u00d79=\"],\n  [\"333\u00d72=\", \"271\u00d72=\"],\n  [\"925\u00d74=\", \"769\u00d75=\"],\n  [\"958\u00d78=\", \"368\u00d77=\"],\n  [\"231\u00d76=\", \"743\u00d76=\"],\n  [\"530\u00d72=\", \"271\u00d75=\"],\n  [\"686\u00d74=\", \"573\u00d74=\"],\n  [\"353\u00d78=\", \"240\u00d78=\"],\n  [\"371\u00d74=\", \"326\u00d78=\"],\n  [\"239\u00d73=\", \"343\u00d79=\"],\n  [\"830\u00d74=\", \"584\u00d78=\"],\n  [\"183\u00d72=\", \"235\u00d79=\"],\n  [\"747\u00d78=\", \"487\u00d78=\"],\n  [\"601\u00d77=\", \"872\u00d72=\"],\n  [\"661\u00d79=\", \"181\u00d76=\"],\n  [\"837\u00d79=\", \"296\u00d79=\"],\n  [\"507\u00d79=\", \"226\u00d74=\"],\n  [\"349\u00d78=\", \"913\u00d78=\"],\n  [\"773\u00d73=\", \"489\u00d72=\"],\n  [\"172\u00d76=\", \"305\u00d74=\"],\n  [\"702\u00d72=\", \"736\u00d72=\"],\n  [\"172\u00d74=\", \"271\u00d73=\"],\n  [\"373\u00d72=\", \"117\u00d76=\"],\n  [\"380\u00d77=\", \"934\u00d73=\"],\n  [\"736\u00d76=\", \"704\u00d75=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-problem cells in the table with new operands.\n# Each \"old\" string is unique in the document, so Find/Replace against the\n# whole document content swaps only the visible text while the existing\n# run formatting (font, size, etc.) stays untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"703\u00d76=\", \"442\u00d79=\"),\n    @(\"333\u00d72=\", \"271\u00d72=\"),\n    @(\"925\u00d74=\", \"769\u00d75=\"),\n    @(\"958\u00d78=\", \"368\u00d77=\"),\n    @(\"231\u00d76=\", \"743\u00d76=\"),\n    @(\"530\u00d72=\", \"271\u00d75=\"),\n    @(\"686\u00d74=\", \"573\u00d74=\"),\n    @(\"353\u00d78=\", \"240\u00d78=\"),\n    @(\"371\u00d74=\", \"326\u00d78=\"),\n    @(\"239\u00d73=\", \"343\u00d79=\"),\n    @(\"830\u00d74=\", \"584\u00d78=\"),\n    @(\"183\u00d72=\", \"235\u00d79=\"),\n    @(\"747\u00d78=\", \"487\u00d78=\"),\n    @(\"601\u00d77=\", \"872\u00d72=\"),\n    @(\"661\u00d79=\", \"181\u00d76=\"),\n    @(\"837\u00d79=\", \"296\u00d79=\"),\n    @(\"507\u00d79=\", \"226\u00d74=\"),\n    @(\"349\u00d78=\", \"913\u00d78=\"),\n    @(\"773\u00d73=\", \"489\u00d72=\"),\n    @(\"172\u00d76=\", \"305\u00d74=\"),\n    @(\"702\u00d72=\", \"736\u00d72=\"),\n    @(\"172\u00d74=\", \"271\u00d73=\"),\n    @(\"373\u00d72=\", \"117\u00d76=\"),\n    @(\"380\u00d77=\", \"934\u00d73=\"),\n    @(\"736\u00d76=\", \"704\u00d75=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
